$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 349116666.8676928
$ws.Range("C2").Value = 405754579.53039044
$ws.Range("D2").Value = 462392492.19308877
$ws.Range("E2").Value = 519030404.85578716
$ws.Range("F2").Value = 575668317.5184852

$ws.Range("B3").Value = 744616858.3258507
$ws.Range("C3").Value = 801254770.9885484
$ws.Range("D3").Value = 857892683.6512467
$ws.Range("E3").Value = 914530596.313945
$ws.Range("F3").Value = 971168508.976643

$ws.Range("B4").Value = 1535970429.3710833
$ws.Range("C4").Value = 1592608342.033781
$ws.Range("D4").Value = 1649246254.6964793
$ws.Range("E4").Value = 1705884167.3591776
$ws.Range("F4").Value = 1762522080.0218756

$ws.Range("B5").Value = 2486216645.587604
$ws.Range("C5").Value = 2542854558.250302
$ws.Range("D5").Value = 2599492470.913
$ws.Range("E5").Value = 2656130383.5756984
$ws.Range("F5").Value = 2712768296.2383966
